# Update XRO ENSO forecasts: revise the latest two forecast rows and
# append a new "2024-10" init row to the Nino34 sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nino34")

# --- small corrections to row 21 (init 2024-08) ---
$ws.Range("I21").Value = -0.408
$ws.Range("P21").Value = -0.638

# --- row 22 (init 2024-09) revised forecast values ---
$ws.Range("B22").Value = -0.217
$ws.Range("C22").Value = -0.295
$ws.Range("D22").Value = -0.427
$ws.Range("E22").Value = -0.498
$ws.Range("F22").Value = -0.482
$ws.Range("G22").Value = -0.419
$ws.Range("H22").Value = -0.38
$ws.Range("I22").Value = -0.397
$ws.Range("J22").Value = -0.476
$ws.Range("K22").Value = -0.597
$ws.Range("L22").Value = -0.6899999999999999
$ws.Range("M22").Value = -0.702
$ws.Range("N22").Value = -0.679
$ws.Range("O22").Value = -0.7
$ws.Range("P22").Value = -0.792
$ws.Range("Q22").Value = -0.897
$ws.Range("R22").Value = -0.922
$ws.Range("S22").Value = -0.839
$ws.Range("T22").Value = -0.7
$ws.Range("U22").Value = -0.555

# --- new row 23 (init 2024-10), matching row 22's look & formatting ---
$label = $ws.Range("A23")
$label.Value = "2024-10"
$label.Font.Bold = $true
$label.Borders.Item(13).LineStyle = 1
$label.HorizontalAlignment = -4108
$label.VerticalAlignment = -4160

$ws.Range("B23").Value = -0.283
$ws.Range("C23").Value = -0.369
$ws.Range("D23").Value = -0.43
$ws.Range("E23").Value = -0.376
$ws.Range("F23").Value = -0.293
$ws.Range("G23").Value = -0.258
$ws.Range("H23").Value = -0.287
$ws.Range("I23").Value = -0.369
$ws.Range("J23").Value = -0.484
$ws.Range("K23").Value = -0.575
$ws.Range("L23").Value = -0.596
$ws.Range("M23").Value = -0.582
$ws.Range("N23").Value = -0.603
$ws.Range("O23").Value = -0.68
$ws.Range("P23").Value = -0.764
$ws.Range("Q23").Value = -0.781
$ws.Range("R23").Value = -0.709
$ws.Range("S23").Value = -0.592
$ws.Range("T23").Value = -0.469
$ws.Range("U23").Value = -0.342

$ws.Range("B23:U23").NumberFormat = $ws.Range("B22:U22").NumberFormat
